# Insert a new weekly price record as row 215 ("Fruta, Vega Monumental
# Concepción - Piña"), pushing the previous rows 215-238 down to 216-239.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 215; existing rows 215..238 shift to 216..239.
$ws.Rows.Item(215).Insert()

# Populate the new row 215 with the new record's data.
$ws.Cells.Item(215, 1).Value = 11
$ws.Cells.Item(215, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(215, 3).Value = "Bíobío"
$ws.Cells.Item(215, 4).Value = 44946
$ws.Cells.Item(215, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215, 5).Value = 8
$ws.Cells.Item(215, 6).Value = "Fruta"
$ws.Cells.Item(215, 7).Value = 100108
$ws.Cells.Item(215, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(215, 9).Value = 100108005
$ws.Cells.Item(215, 10).Value = "Piña"
$ws.Cells.Item(215, 11).Value = "Caramelo"
$ws.Cells.Item(215, 12).Value = "Primera"
$ws.Cells.Item(215, 13).Value = 200
$ws.Cells.Item(215, 14).Value = 18000
$ws.Cells.Item(215, 15).Value = 19000
$ws.Cells.Item(215, 16).Value = 18500
$ws.Cells.Item(215, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(215, 18).Value = "Ecuador"
$ws.Cells.Item(215, 19).Value = 1542
$ws.Cells.Item(215, 20).Value = 12
